# Finished Week 13 logging
# Update target depth data on both the OFF and DEF sheets, row 2 (Home row)

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 564
$wsOff.Range("C2").Value = 409
$wsOff.Range("D2").Value = 110
$wsOff.Range("E2").Value = 46
$wsOff.Range("F2").Value = 8

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 352
$wsDef.Range("C2").Value = 219
$wsDef.Range("D2").Value = 83
$wsDef.Range("E2").Value = 30
$wsDef.Range("F2").Value = 9
